$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Category rows (3-12): new "Tên yêu cầu" values ---
$ws.Cells.Item(3,2).Value  = "Quản lí nhân viên"
$ws.Cells.Item(4,2).Value  = "Chấm công"
$ws.Cells.Item(5,2).Value  = "Quản lí danh mục"
$ws.Cells.Item(6,2).Value  = "Báo cáo"
$ws.Cells.Item(7,2).Value  = "Tìm kiếm"
$ws.Cells.Item(8,2).Value  = "Quản lí kho"
$ws.Cells.Item(9,2).Value  = "Quản lí đặt hàng"
$ws.Cells.Item(10,2).Value = "Quản lí hóa đơn"
$ws.Cells.Item(11,2).Value = "Quản lí khách hàng"
$ws.Cells.Item(12,2).Value = "Quản lí hệ thống"

# New related-object note for the "Chấm công" row
$ws.Cells.Item(4,3).Value = "Máy chấm công bằng vân tay"

# --- Clear obsolete "Tên yêu cầu" bullet list (rows 13-36) ---
for ($r = 13; $r -le 36; $r++) {
    $ws.Cells.Item($r,2).Value = ""
}

# --- Title row: update last so it becomes the last new shared string ---
$ws.Cells.Item(1,1).Value = "Yêu cầu tương thích"

# --- Restore view: scroll back to top-left and select the title row ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1:D1").Select()
